# Regenerate save_data column G ("K") values (was based on Strike#, now
# recalculated/regenerated). Update each row's K value in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 4
    4  = 4
    5  = 6
    6  = 6
    7  = 3
    8  = 5
    9  = 5
    10 = 5
    11 = 8
    12 = 0
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    17 = 3
    18 = 1
    19 = 3
    20 = 0
    21 = 4
    22 = 1
    23 = 0
    24 = 1
    25 = 2
    26 = 3
    27 = 2
    28 = 1
    29 = 0
    30 = 1
    31 = 2
    32 = 1
    33 = 0
    34 = 2
    35 = 2
    36 = 1
    37 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
